$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 10974.25
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

$ws.Range("H19").Value = 1301.2106
$ws.Range("J19").Value = 1763.7273
$ws.Range("L19").Value = 1763.7273
$ws.Range("N19").Value = -2113.7273

$ws.Range("H41").Value = 1450.6666
$ws.Range("I41").Value = 1384.7778
$ws.Range("J41").Value = 1549.5
$ws.Range("K41").Value = 1384.7778
$ws.Range("L41").Value = 1549.5
$ws.Range("M41").Value = -944.7778000000001
$ws.Range("N41").Value = -2429.5

$ws.Range("H100").Value = 1256.8
$ws.Range("I100").Value = 1071
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1071
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -530
$ws.Range("N100").Value = -3082

$ws.Range("H111").Value = 4836.2856
$ws.Range("I111").Value = 4975.8335
$ws.Range("K111").Value = 14927.5005
$ws.Range("M111").Value = -11860.5005

$ws.Range("H141").Value = 3258.182
$ws.Range("I141").Value = 3315.5557
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 9946.667099999999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -4766.667099999999
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9436877
$ws.Range("J32").Value = 13982.667
$ws.Range("L32").Value = 13982.667
$ws.Range("N32").Value = -14556.667

$ws.Range("H45").Value = 2149
$ws.Range("I45").Value = 1342.875
$ws.Range("K45").Value = 1342.875
$ws.Range("M45").Value = -965.875

$ws.Range("H46").Value = 4901
$ws.Range("J46").Value = 4901
$ws.Range("L46").Value = 4901
$ws.Range("N46").Value = -5539

$ws.Range("H63").Value = 4998
$ws.Range("I63").Value = 4331
$ws.Range("K63").Value = 4331
$ws.Range("M63").Value = -3645

$ws.Range("H66").Value = 4998
$ws.Range("I66").Value = 4331
$ws.Range("K66").Value = 21655
$ws.Range("M66").Value = -18223

$ws.Range("H104").Value = 48130
$ws.Range("J104").Value = 48130
$ws.Range("L104").Value = 48130
$ws.Range("N104").Value = -55118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 823.6429000000001
$ws.Range("J80").Value = 682.4
$ws.Range("L80").Value = 682.4
$ws.Range("N80").Value = -2678.4

$ws.Range("H83").Value = 823.6429000000001
$ws.Range("J83").Value = 682.4
$ws.Range("L83").Value = 3412
$ws.Range("N83").Value = -13396

$ws.Range("H105").Value = 2938.1428
$ws.Range("I105").Value = 2266.75
$ws.Range("K105").Value = 2266.75
$ws.Range("M105").Value = -519.75

$ws.Range("H134").Value = 48148.5
$ws.Range("I134").Value = 802.7895
$ws.Range("J134").Value = 348004.66
$ws.Range("K134").Value = 2408.3685
$ws.Range("L134").Value = 1044013.98
$ws.Range("M134").Value = 126.6315
$ws.Range("N134").Value = -1049083.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = 0

$ws.Range("H31").Value = 637723.3
$ws.Range("I31").Value = 9315.259
$ws.Range("K31").Value = 9315.259
$ws.Range("M31").Value = -9020.259

$ws.Range("H34").Value = 637723.3
$ws.Range("I34").Value = 9315.259
$ws.Range("K34").Value = 9315.259
$ws.Range("M34").Value = -9113.259

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""

$ws.Range("H58").Value = 1966.3077
$ws.Range("I58").Value = 2058
$ws.Range("J58").Value = 1760
$ws.Range("K58").Value = 2058
$ws.Range("L58").Value = 1760
$ws.Range("M58").Value = -1855
$ws.Range("N58").Value = -2166

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = 0

$ws.Range("H136").Value = 1966.3077
$ws.Range("I136").Value = 2058
$ws.Range("J136").Value = 1760
$ws.Range("K136").Value = 6174
$ws.Range("L136").Value = 5280
$ws.Range("M136").Value = -3624
$ws.Range("N136").Value = -10380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4539577
$ws.Range("I4").Value = 15132062
$ws.Range("K4").Value = 45396186
$ws.Range("M4").Value = -45396074

$ws.Range("H11").Value = 2139.5
$ws.Range("I11").Value = 2169.08
$ws.Range("K11").Value = 6507.24
$ws.Range("M11").Value = -6367.24

$ws.Range("H92").Value = 1430302.8
$ws.Range("J92").Value = 417.25
$ws.Range("L92").Value = 1251.75
$ws.Range("N92").Value = -3747.75

$ws.Range("H133").Value = 5701.846
$ws.Range("I133").Value = 5124.8887
$ws.Range("K133").Value = 15374.6661
$ws.Range("M133").Value = -10314.6661

$ws.Range("H137").Value = 7172.1665
$ws.Range("J137").Value = 5758.25
$ws.Range("L137").Value = 17274.75
$ws.Range("N137").Value = -27474.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1762.6086
$ws.Range("I107").Value = 1147.9375
$ws.Range("J107").Value = 3167.5715
$ws.Range("K107").Value = 1147.9375
$ws.Range("L107").Value = 3167.5715
$ws.Range("M107").Value = 772.0625
$ws.Range("N107").Value = -7007.5715

$ws.Range("H126").Value = 4008.2
$ws.Range("I126").Value = 3918.2
$ws.Range("J126").Value = 4098.2
$ws.Range("K126").Value = 11754.6
$ws.Range("L126").Value = 12294.6
$ws.Range("M126").Value = -9284.599999999999
$ws.Range("N126").Value = -17234.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6001
$ws.Range("I22").Value = 6001
$ws.Range("K22").Value = 6001
$ws.Range("M22").Value = -5706

$ws.Range("H27").Value = 6001
$ws.Range("I27").Value = 6001
$ws.Range("K27").Value = 6001
$ws.Range("M27").Value = -5894

$ws.Range("H40").Value = 3606.8708
$ws.Range("I40").Value = 3050.5833
$ws.Range("J40").Value = 5514.143
$ws.Range("K40").Value = 3050.5833
$ws.Range("L40").Value = 5514.143
$ws.Range("M40").Value = -2914.5833
$ws.Range("N40").Value = -5786.143

$ws.Range("H68").Value = 2999.5
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250

$ws.Range("H71").Value = 2999.5
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251

$ws.Range("H100").Value = 1284.5714
$ws.Range("I100").Value = 1284.5714
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1284.5714
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -743.5714

$ws.Range("H122").Value = 5361
$ws.Range("I122").Value = 4633.1875
$ws.Range("J122").Value = 7690
$ws.Range("K122").Value = 13899.5625
$ws.Range("L122").Value = 23070
$ws.Range("M122").Value = -11449.5625
$ws.Range("N122").Value = -27970

$ws.Range("H127").Value = 115828.75
$ws.Range("J127").Value = 115828.75
$ws.Range("L127").Value = 115828.75
$ws.Range("N127").Value = -125748.75

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

$ws.Range("H136").Value = 38756.83
$ws.Range("I136").Value = 6316.3335
$ws.Range("J136").Value = 73105.586
$ws.Range("K136").Value = 18949.0005
$ws.Range("L136").Value = 219316.758
$ws.Range("M136").Value = -16399.0005
$ws.Range("N136").Value = -224416.758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5010
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""

$ws.Range("H41").Value = 21657.334
$ws.Range("J41").Value = 21657.334
$ws.Range("L41").Value = 21657.334
$ws.Range("N41").Value = -22437.334

$ws.Range("H107").Value = 45459824
$ws.Range("I107").Value = 62502708
$ws.Range("J107").Value = 12133.667
$ws.Range("K107").Value = 187508124
$ws.Range("L107").Value = 36401.001
$ws.Range("M107").Value = -187506204
$ws.Range("N107").Value = -40241.001

$ws.Range("H110").Value = 34000
$ws.Range("J110").Value = 34000
$ws.Range("L110").Value = 34000
$ws.Range("N110").Value = -42180
